$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, shifting existing rows 62-72 down to 63-73
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new data record
$ws.Cells.Item(62,1).Value = 9
$ws.Cells.Item(62,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(62,3).Value = "Metropolitana"
$ws.Cells.Item(62,4).Value = 44504
$ws.Cells.Item(62,5).Value = 13
$ws.Cells.Item(62,6).Value = 100112022
$ws.Cells.Item(62,7).Value = "Arveja Verde"
$ws.Cells.Item(62,8).Value = "Sin especificar"
$ws.Cells.Item(62,9).Value = "Primera"
$ws.Cells.Item(62,10).Value = 26
$ws.Cells.Item(62,11).Value = 12000
$ws.Cells.Item(62,12).Value = 13000
$ws.Cells.Item(62,13).Value = 12500
$ws.Cells.Item(62,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(62,15).Value = "Región Metropolitana"
$ws.Cells.Item(62,16).Value = 500
$ws.Cells.Item(62,17).Value = 25
$ws.Cells.Item(62,18).Value = "Hortaliza"
